$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 12
$ws.Range("F2").Value = 4
$ws.Range("G2").Value = 0

$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 19
$ws.Range("F3").Value = 7
$ws.Range("G3").Value = 0

$ws.Range("A4").Value = 0
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 5
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 20
$ws.Range("F4").Value = 4
$ws.Range("G4").Value = 0
